$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the Price (D) / Volume(1h) (E) columns with the latest scraped
# figures. Both columns hold plain TEXT in this workbook (e.g. "215.61",
# "0.5230", "26.060.53", "  -0.25%  ") rather than numbers. A handful of
# the new Price values look like ordinary decimals (e.g. "0.5230",
# "215.77") and would otherwise be auto-detected by Excel's normal
# Value-assignment as numbers -- silently dropping significant trailing
# zeros ("0.5230" -> 0.523) and changing the cell's stored type. Those
# cells are written via Formula with a leading apostrophe (Excel's
# standard "force text" quote-prefix) so they stay literal text, exactly
# like every other cell in these two columns; the rest use plain Value
# assignment since they're already unambiguous as text (percentages,
# the "thousands-dot" prices, or values with the subscript-zero glyph).
$ws.Range("D2").Value = "26.069.47"
$ws.Range("E2").Value = "  -0.20%  "
$ws.Range("D3").Value = "1.647.07"
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").Formula = "'215.77"
$ws.Range("E5").Value = "  +2.30%  "
$ws.Range("D6").Formula = "'0.5230"
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("D8").Formula = "'0.2614"
$ws.Range("E8").Value = "  -0.41%  "
$ws.Range("D9").Formula = "'0.06366"
$ws.Range("E9").Value = "  +0.50%  "
$ws.Range("D10").Formula = "'20.89"
$ws.Range("E10").Value = "  -1.56%  "
$ws.Range("D11").Formula = "'0.07655"
$ws.Range("E11").Value = "  +1.34%  "
$ws.Range("D12").Value = "1.647.48"
$ws.Range("E12").Value = "  -1.38%  "
$ws.Range("D13").Formula = "'4.426"
$ws.Range("E13").Value = "  -0.15%  "
$ws.Range("D14").Value = "1.867.83"
$ws.Range("E14").Value = "  -1.61%  "
$ws.Range("D15").Formula = "'0.5553"
$ws.Range("E15").Value = "  +1.21%  "
$ws.Range("D16").Value = "0.0₅8300"
$ws.Range("E16").Value = "  +3.12%  "
$ws.Range("E17").Value = "  -2.24%  "
$ws.Range("D18").Value = "26.073.02"
$ws.Range("E18").Value = "  -0.40%  "
$ws.Range("E19").Value = "  -0.20%  "
$ws.Range("E20").Value = "  -0.54%  "
$ws.Range("D21").Formula = "'188.46"
$ws.Range("E21").Value = "  +0.27%  "
$ws.Range("E22").Value = "  -1.02%  "
$ws.Range("D23").Formula = "'6.256"
$ws.Range("E23").Value = "  +0.30%  "
$ws.Range("E24").Value = "  -0.20%  "
$ws.Range("D25").Formula = "'146.03"
$ws.Range("E25").Value = "  -2.35%  "
$ws.Range("D26").Formula = "'0.1220"
$ws.Range("E26").Value = "  -1.72%  "
$ws.Range("D27").Formula = "'7.416"
$ws.Range("E27").Value = "  -0.85%  "
$ws.Range("D28").Formula = "'15.83"
$ws.Range("E28").Value = "  +0.12%  "
$ws.Range("D29").Formula = "'1.396"
$ws.Range("E29").Value = "  +3.02%  "
$ws.Range("D30").Formula = "'0.05958"
$ws.Range("E30").Value = "  -5.94%  "
$ws.Range("D31").Formula = "'1.268"
$ws.Range("E31").Value = "  -1.16%  "
$ws.Range("D32").Formula = "'3.402"
$ws.Range("E32").Value = "  -0.42%  "
$ws.Range("D33").Formula = "'3.401"
$ws.Range("E33").Value = "  -3.61%  "
$ws.Range("E34").Value = "  +0.50%  "
$ws.Range("D35").Formula = "'0.9981"
$ws.Range("E35").Value = "  -0.68%  "
$ws.Range("D36").Formula = "'2.393"
$ws.Range("E36").Value = "  -0.10%  "
$ws.Range("D37").Formula = "'2.755"
$ws.Range("E37").Value = "  -0.28%  "
$ws.Range("D38").Formula = "'0.5638"
$ws.Range("E38").Value = "  -6.54%  "
$ws.Range("D39").Formula = "'0.01609"
$ws.Range("E39").Value = "  -0.21%  "
$ws.Range("D40").Formula = "'0.8567"
$ws.Range("E40").Value = "  -0.88%  "
$ws.Range("D41").Formula = "'5.838"
$ws.Range("E41").Value = "  -3.73%  "
$ws.Range("E42").Value = "  -0.24%  "
$ws.Range("D43").Value = "1.026.93"
$ws.Range("E43").Value = "  -8.15%  "
$ws.Range("D44").Formula = "'99.26"
$ws.Range("E44").Value = "  -1.20%  "
$ws.Range("D45").Value = "1.795.49"
$ws.Range("E45").Value = "  -1.40%  "
$ws.Range("D46").Value = "0.0₈110"
$ws.Range("E46").Value = "  +1.74%  "
$ws.Range("D47").Formula = "'55.78"
$ws.Range("E47").Value = "  +0.58%  "
$ws.Range("E48").Value = "  +0.13%  "
$ws.Range("D49").Formula = "'8.085"
$ws.Range("E49").Value = "  +0.27%  "
$ws.Range("E50").Value = "  -1.66%  "
$ws.Range("D51").Formula = "'0.4216"
$ws.Range("E51").Value = "  -0.55%  "
